$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Framework")

$states = @(
    "Abia",
    "Adamawa",
    "Akwa lbom",
    "Anambra",
    "Bauchi",
    "Bayelsa",
    "Benue",
    "Borno",
    "Cross River",
    "Delta",
    "Ebonyi",
    "Edo",
    "Ekiti",
    "Enugu",
    "Federal Capital Territory",
    "Gombe",
    "Imo",
    "Jigawa",
    "Kaduna",
    "Kano",
    "Katsina",
    "Kebbi",
    "Kogi",
    "Kwara",
    "Lagos",
    "Nasarawa",
    "Niger",
    "Ogun",
    "Ondo",
    "Osun",
    "Oyo",
    "Plateau",
    "Rivers",
    "Sokoto",
    "Taraba",
    "Yobe"
)

$col = 8
foreach ($state in $states) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $state
    $cell.Style = "Normal"
    if ($col -eq 8) {
        $cell.Font.Bold = $true
    }
    $col = $col + 1
}

$ws.Columns.Item(10).AutoFit() | Out-Null
$ws.Columns.Item(16).AutoFit() | Out-Null
$ws.Columns.Item(22).AutoFit() | Out-Null
